$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.035.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.832.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.50%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6260'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.56%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07579'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2921'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.11%  '

$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.830.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.957'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6642'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.34%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001009'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +15.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.026'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.025.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.61%  '

$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.079.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.47%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.168'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.35%  '

$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.55%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.487'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.56%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1374'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.36%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.65%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.491'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.55%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.097'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.01%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.009'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.192'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.55%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05204'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.36%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.843'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.91%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7350'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.62%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.139'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.65%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.697'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.73%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.239.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.66%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.758'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01783'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.76%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.325'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.52%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8959'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.12%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.23%  '

$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.978.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.62%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000124'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.20%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.55%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5108'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.65%  '

$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4025'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.876'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.20%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05755'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.94%  '
